$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the split file info
$ws.Range("B2").Value = "E:\storage\128G.mp4"
$ws.Range("C2").Value = "ád"
$ws.Range("D2").Value = "ád"
$ws.Range("F2").Value = "17:52"
$ws.Range("G2").Value = "E:/New folder\128G.mp4"

# Remove rows 3 through 6 (the rest of the split parts)
$ws.Range("A3:G6").EntireRow.Delete()
